$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PM1_Pu_GHS")

# Rows 15-28 had their regression-parameter scratch columns (AL, AR:AZ, BB,
# BC, BE, BF, BH) cleared out - only rows 1-14 keep that data now.
$ws.Range("AL15:AL28").ClearContents()
$ws.Range("AR15:AZ28").ClearContents()
$ws.Range("BB15:BC28").ClearContents()
$ws.Range("BE15:BF28").ClearContents()
$ws.Range("BH15:BH28").ClearContents()
